$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (not auto-converted to numbers).
# Apply a Text number format to the whole data range first, write values, then reset cell
# style back to Normal (index 0) so no stray style index is left on the cells, matching the
# original (unstyled) cells exactly.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("B44").Value = 'mCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'

# Column D (price) updates
$ws.Range("D2").Value = '27.456.18'
$ws.Range("D3").Value = '1.616.45'
$ws.Range("D5").Value = '211.02'
$ws.Range("D8").Value = '22.74'
$ws.Range("D10").Value = '0.0611'
$ws.Range("D12").Value = '1.847.41'
$ws.Range("D13").Value = '1.627.78'
$ws.Range("D15").Value = '0.548'
$ws.Range("D16").Value = '65.10'
$ws.Range("D17").Value = '27.436.23'
$ws.Range("D18").Value = '232.62'
$ws.Range("D19").Value = '0.0₃0716'
$ws.Range("D20").Value = '7.52'
$ws.Range("D23").Value = '10.20'
$ws.Range("D25").Value = '150.41'
$ws.Range("D26").Value = '6.85'
$ws.Range("D29").Value = '15.52'
$ws.Range("D31").Value = '0.0482'
$ws.Range("D33").Value = '1.469.28'
$ws.Range("D34").Value = '3.06'
$ws.Range("D37").Value = '0.963'
$ws.Range("D39").Value = '0.556'
$ws.Range("D40").Value = '0.860'
$ws.Range("D42").Value = '67.02'
$ws.Range("D43").Value = '0.986'
$ws.Range("D44").Value = '2.46'
$ws.Range("D47").Value = '1.757.81'
$ws.Range("D48").Value = '1.72'
$ws.Range("D49").Value = '86.80'
$ws.Range("D50").Value = '0.0₆0104'

# Reset column D style back to Normal so the text-format tweak above does not leave a
# lingering style index on the data cells (matches original, style-less cells).
$dRange.Style = "Normal"

# Column E (1h volume/change) updates
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("E8").Value = '  -1.99%  '
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("E15").Value = '  -2.97%  '
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  +5.32%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  -1.92%  '
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("E34").Value = '  -3.38%  '
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("E37").Value = '  +9.49%  '
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("E40").Value = '  -3.07%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("E43").Value = '  -4.90%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  -2.69%  '
$ws.Range("E46").Value = '  -6.74%  '
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("E50").Value = '  -1.94%  '
$ws.Range("E51").Value = '  +0.72%  '
